# Weekly fruit/vegetable price update:
# Insert a new daily record at the top of the data block (row 23),
# pushing all existing records down by one row (row 61 -> row 62).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 23:61 down by inserting a new blank row at 23.
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row with the latest market record.
$ws.Range("A23").Value = 9
$ws.Range("B23").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C23").Value = "Metropolitana"
$ws.Range("D23").Value = 44797
$ws.Range("E23").Value = 13
$ws.Range("F23").Value = 100112035
$ws.Range("G23").Value = "Bruselas (repollito)"
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 43
$ws.Range("K23").Value = 20000
$ws.Range("L23").Value = 20000
$ws.Range("M23").Value = 20000
$ws.Range("N23").Value = "`$/malla 15 kilos"
$ws.Range("O23").Value = "Hijuelas"
$ws.Range("P23").Value = 1333
$ws.Range("Q23").Value = 15
$ws.Range("R23").Value = "Hortaliza"
